$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that looks numeric (e.g. "9.00", "13.90").
# Force it to Text format first so Excel keeps the exact string (trailing
# zeros, thousand-separator dots, etc.) instead of silently coercing the
# assigned value to a float and dropping formatting.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "35.445.71"
$ws.Range("E2").Value = "  -3.65%  "
$ws.Range("D3").Value = "1.975.08"
$ws.Range("E3").Value = "  -5.56%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "239.61"
$ws.Range("E5").Value = "  -2.24%  "
$ws.Range("E6").Value = "  -4.17%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "55.27"
$ws.Range("E8").Value = "  +2.10%  "
$ws.Range("D9").Value = "58.94"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "0.351"
$ws.Range("E10").Value = "  -4.31%  "
$ws.Range("E11").Value = "  -5.68%  "
$ws.Range("D12").Value = "0.103"
$ws.Range("E12").Value = "  -5.69%  "
$ws.Range("D13").Value = "0.875"
$ws.Range("E13").Value = "  -3.75%  "
$ws.Range("D14").Value = "13.90"
$ws.Range("E14").Value = "  -7.56%  "
$ws.Range("D15").Value = "2.269.52"
$ws.Range("E15").Value = "  -5.35%  "
$ws.Range("E16").Value = "  -6.22%  "
$ws.Range("D17").Value = "1.984.04"
$ws.Range("E17").Value = "  -5.30%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").Value = "16.95"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "35.337.92"
$ws.Range("E19").Value = "  -3.91%  "
$ws.Range("E20").Value = "  -4.31%  "
$ws.Range("D21").Value = "0.0₃0824"
$ws.Range("E21").Value = "  -6.23%  "
$ws.Range("D22").Value = "230.02"
$ws.Range("E22").Value = "  -3.60%  "
$ws.Range("D23").Value = "4.92"
$ws.Range("E23").Value = "  -9.80%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "2.31"
$ws.Range("E25").Value = "  -3.69%  "
$ws.Range("D26").Value = "2.21"
$ws.Range("E26").Value = "  +2.66%  "
$ws.Range("D27").Value = "162.08"
$ws.Range("E27").Value = "  -2.72%  "
$ws.Range("D28").Value = "9.00"
$ws.Range("E28").Value = "  -7.53%  "
$ws.Range("D29").Value = "19.20"
$ws.Range("E29").Value = "  -8.20%  "
$ws.Range("E30").Value = "  -3.91%  "
$ws.Range("E31").Value = "  -3.52%  "
$ws.Range("D32").Value = "4.69"
$ws.Range("E32").Value = "  -10.41%  "
$ws.Range("E33").Value = "  -4.70%  "
$ws.Range("D34").Value = "0.0883"
$ws.Range("E34").Value = "  +6.99%  "
$ws.Range("D35").Value = "4.19"
$ws.Range("E35").Value = "  -10.21%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  -2.19%  "
$ws.Range("D38").Value = "2.16"
$ws.Range("E38").Value = "  -10.96%  "
$ws.Range("E39").Value = "  -2.02%  "
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("E41").Value = "  -8.67%  "
$ws.Range("E42").Value = "  -6.64%  "
$ws.Range("E43").Value = "  -7.87%  "
$ws.Range("D44").Value = "1.361.52"
$ws.Range("E44").Value = "  -2.95%  "
$ws.Range("D45").Value = "0.0867"
$ws.Range("E45").Value = "  -9.01%  "
$ws.Range("D46").Value = "88.59"
$ws.Range("E46").Value = "  -7.92%  "
$ws.Range("D47").Value = "7.29"
$ws.Range("E47").Value = "  -2.75%  "
$ws.Range("D48").Value = "15.18"
$ws.Range("E48").Value = "  -5.75%  "
$ws.Range("D49").Value = "2.88"
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("D50").Value = "2.21"
$ws.Range("E50").Value = "  -10.27%  "
$ws.Range("E51").Value = "  -2.17%  "
